$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Tiny Toy Tanks"
$ws.Range("A3").Value = "Super Buckyball Tournament Preseason"
$ws.Range("A4").Value = "Beat Me! - Puppetonia Tournament"
$ws.Range("A5").Value = "Perfect Vermin"
